$d = $word.ActiveDocument

# --- change1 (merge 'Ten eerste is het mogelijk ' + 'om ') ---
$xml1 = @'
<w:p w14:paraId="0E1B82DE" w14:textId="73FC4862" w:rsidR="007D2394" w:rsidRPr="00B02C1C" w:rsidRDefault="007D2394" w:rsidP="007D2394"><w:r><w:t xml:space="preserve">Ten eerste is het mogelijk om </w:t></w:r><w:r w:rsidR="00903EBD" w:rsidRPr="00B02C1C"><w:t xml:space="preserve">“Native” applicaties te ontwikkelen. Hiermee worden de applicaties bedoeld die zijn ontwikkeld met de software die de makers van de twee besturingssystemen hebben meegeleverd voor de ontwikkeling van de applicaties.  </w:t></w:r></w:p>
'@
$found1 = $false
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $para = $d.Paragraphs($i)
    if ($para.Range.Text.StartsWith("Ten eerste is het mogelijk")) {
        $null = $para.Range.InsertXML($xml1)
        $found1 = $true
        break
    }
}
Write-Output "change1 found=$found1"

# --- change2 (merge Hybride applicaties runs) ---
$xml2 = @'
<w:p w14:paraId="67B6610B" w14:textId="40BCD597" w:rsidR="00F36006" w:rsidRPr="00F36006" w:rsidRDefault="00F36006" w:rsidP="007913AD"><w:r><w:t>Hybride applicaties zijn zoals crossplatform applicaties geschreven met een code base die werkt op verschillende besturingssystemen en apparaten. Het voornaamste verschil tussen hybride en crossplatform is de opbouw van de gebruikersinterface. Waarbij de gebruikersinterface van hybride apps is geschreven met web technologieën.</w:t></w:r></w:p>
'@
$found2 = $false
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $para = $d.Paragraphs($i)
    if ($para.Range.Text.StartsWith("Hybride applicaties zijn")) {
        $null = $para.Range.InsertXML($xml2)
        $found2 = $true
        break
    }
}
Write-Output "change2 found=$found2"

# --- change3 (fix 'Hierbij worden' double space + remove proofErr) ---
$xml3 = @'
<w:p w14:paraId="7F098F3E" w14:textId="13E5BF20" w:rsidR="00B02C1C" w:rsidRPr="00B02C1C" w:rsidRDefault="00B02C1C" w:rsidP="007913AD"><w:r w:rsidRPr="00B02C1C"><w:t>“</w:t></w:r><w:proofErr w:type="spellStart"/><w:r w:rsidRPr="00B02C1C"><w:t>Progressive</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r w:rsidRPr="00B02C1C"><w:t xml:space="preserve"> web apps” zijn webapplicaties ontwikkeld om te optimaal te kunnen worden gebruikt op mobiele apparatuur. De applicatie wordt altijd geopend door een webbrowser zonder dat daarvoor installatie nodig is.  </w:t></w:r><w:r><w:t xml:space="preserve">Hierbij worden geen afzonderlijke versies voor de verschillende besturingssystemen gebruikt. </w:t></w:r></w:p>
'@
$found3 = $false
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $para = $d.Paragraphs($i)
    if ($para.Range.Text.StartsWith("“Progressive")) {
        $null = $para.Range.InsertXML($xml3)
        $found3 = $true
        break
    }
}
Write-Output "change3 found=$found3"

# --- change4 (de succes -> het succes) ---
$xml4 = @'
<w:p w14:paraId="4C60A452" w14:textId="60694CDE" w:rsidR="00C06906" w:rsidRDefault="00C06906" w:rsidP="00C06906"><w:pPr><w:pStyle w:val="Lijstalinea"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="2"/></w:numPr></w:pPr><w:r><w:t>Gebruikerservaring: De gebruikerservaring is een belangrijke factor voor</w:t></w:r><w:r><w:t xml:space="preserve"> het</w:t></w:r><w:r><w:t xml:space="preserve"> </w:t></w:r><w:r w:rsidR="006E2655"><w:t>succes</w:t></w:r><w:r><w:t xml:space="preserve"> van een applicatie. </w:t></w:r><w:r w:rsidR="006E2655"><w:t xml:space="preserve">Native applicaties kunnen de beste gebruikerservaring bieden, terwijl native en hybride applicaties wellicht concessies moeten maken op verschillende aspecten. </w:t></w:r></w:p>
'@
$found4 = $false
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $para = $d.Paragraphs($i)
    if ($para.Range.Text.StartsWith("Gebruikerservaring: De gebruikerservaring")) {
        $null = $para.Range.InsertXML($xml4)
        $found4 = $true
        break
    }
}
Write-Output "change4 found=$found4"

